# Updates per scheduled-runner market-data refresh (see commit message:
# "chore: update Sheets via scheduled runner").
# For each affected Leve row, refresh the market-price columns
# (H..N: currentAveragePrice / *NQ / *HQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 3713.3333
$ws.Range("J17").Value = 5070
$ws.Range("L17").Value = 15210
$ws.Range("N17").Value = -15546

# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 303
$ws.Range("I33").Value = 151.1
$ws.Range("K33").Value = 151.1
$ws.Range("M33").Value = 77.90000000000001

# Row 125: Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 1618.1666
$ws.Range("I125").Value = 1317.1428
$ws.Range("K125").Value = 11854.2852
$ws.Range("M125").Value = -9394.2852

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 6805448.5
$ws.Range("I138").Value = 1417.5
$ws.Range("J138").Value = 8133064
$ws.Range("K138").Value = 4252.5
$ws.Range("L138").Value = 24399192
$ws.Range("M138").Value = 887.5
$ws.Range("N138").Value = -24409472

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 4889.8
$ws.Range("I45").Value = 2987.25
$ws.Range("K45").Value = 2987.25
$ws.Range("M45").Value = -2610.25

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 6551.778
$ws.Range("I61").Value = 5593.4
$ws.Range("K61").Value = 5593.4
$ws.Range("M61").Value = -5381.4

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 1607.4615
$ws.Range("I74").Value = 1607.4615
$ws.Range("K74").Value = 1607.4615
$ws.Range("M74").Value = -733.4614999999999

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 1607.4615
$ws.Range("I77").Value = 1607.4615
$ws.Range("K77").Value = 8037.307499999999
$ws.Range("M77").Value = -3669.307499999999

# Row 92: Mail It In / High Steel Scale Mail of Fending
$ws.Range("H92").Value = 97887
$ws.Range("J92").Value = 97887
$ws.Range("L92").Value = 97887
$ws.Range("N92").Value = -102879

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 4119.606
$ws.Range("I122").Value = 2449.0322
$ws.Range("J122").Value = 30013.5
$ws.Range("K122").Value = 7347.096600000001
$ws.Range("L122").Value = 90040.5
$ws.Range("M122").Value = -4897.096600000001
$ws.Range("N122").Value = -94940.5

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 6937.9414
$ws.Range("I132").Value = 7213.273
$ws.Range("K132").Value = 21639.819
$ws.Range("M132").Value = -19109.819

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 6551.778
$ws.Range("I136").Value = 5593.4
$ws.Range("K136").Value = 16780.2
$ws.Range("M136").Value = -14230.2

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 5582.316
$ws.Range("I86").Value = 3905.1667
$ws.Range("J86").Value = 8457.429
$ws.Range("K86").Value = 3905.1667
$ws.Range("L86").Value = 8457.429
$ws.Range("M86").Value = -2782.1667
$ws.Range("N86").Value = -10703.429

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 5582.316
$ws.Range("I89").Value = 3905.1667
$ws.Range("J89").Value = 8457.429
$ws.Range("K89").Value = 19525.8335
$ws.Range("L89").Value = 42287.145
$ws.Range("M89").Value = -13909.8335
$ws.Range("N89").Value = -53519.145

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 2128.75
$ws.Range("J58").Value = 5375
$ws.Range("L58").Value = 5375
$ws.Range("N58").Value = -5781

# Row 105: Zelkova, My Love / Zelkova Lumber
$ws.Range("H105").Value = 622.56665
$ws.Range("I105").Value = 452.95
$ws.Range("K105").Value = 452.95
$ws.Range("M105").Value = 1294.05

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 6318504
$ws.Range("I122").Value = 9827662
$ws.Range("J122").Value = 2019.9
$ws.Range("K122").Value = 29482986
$ws.Range("L122").Value = 6059.700000000001
$ws.Range("M122").Value = -29480536
$ws.Range("N122").Value = -10959.7

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 3275.9773
$ws.Range("I132").Value = 2913.5642
$ws.Range("K132").Value = 8740.692599999998
$ws.Range("M132").Value = -6210.692599999998

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 4637.125
$ws.Range("I134").Value = 3019.6
$ws.Range("K134").Value = 9058.799999999999
$ws.Range("M134").Value = -6523.799999999999

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 2128.75
$ws.Range("J136").Value = 5375
$ws.Range("L136").Value = 16125
$ws.Range("N136").Value = -21225

$ws = $wb.Worksheets.Item("CUL")
# Row 97: The Frier Never Lies / Cottonseed Oil
$ws.Range("H97").Value = 1083.6154
$ws.Range("I97").Value = 897.25
$ws.Range("J97").Value = 1166.4445
$ws.Range("K97").Value = 2691.75
$ws.Range("L97").Value = 3499.3335
$ws.Range("M97").Value = -2195.75
$ws.Range("N97").Value = -4491.333500000001

# Row 116: On a Full Stomach / Sausage Links
$ws.Range("H116").Value = 1576.5555
$ws.Range("I116").Value = 998.3333
$ws.Range("K116").Value = 2994.9999
$ws.Range("M116").Value = 447.0001000000002

# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 1948.1666
$ws.Range("J122").Value = 2012.9231
$ws.Range("L122").Value = 18116.3079
$ws.Range("N122").Value = -23016.3079

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 19676460
$ws.Range("I80").Value = 114248.9
$ws.Range("J80").Value = 47622476
$ws.Range("K80").Value = 114248.9
$ws.Range("L80").Value = 47622476
$ws.Range("M80").Value = -113250.9
$ws.Range("N80").Value = -47624472

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 19676460
$ws.Range("I83").Value = 114248.9
$ws.Range("J83").Value = 47622476
$ws.Range("K83").Value = 571244.5
$ws.Range("L83").Value = 238112380
$ws.Range("M83").Value = -566252.5
$ws.Range("N83").Value = -238122364

$ws = $wb.Worksheets.Item("LTW")
# Row 4: Sole Traders / Leather Duckbills
$ws.Range("H4").Value = 750
$ws.Range("I4").Value = 750
$ws.Range("K4").Value = 750
$ws.Range("M4").Value = -637

# Row 28: My Sole to Take / Padded Leather Duckbills
$ws.Range("H28").Value = 750
$ws.Range("I28").Value = 750
$ws.Range("K28").Value = 750
$ws.Range("M28").Value = -518

# Row 37: Quicker than Sand / Padded Leather Duckbills
$ws.Range("H37").Value = 750
$ws.Range("I37").Value = 750
$ws.Range("K37").Value = 750
$ws.Range("M37").Value = -643

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 2519.195
$ws.Range("I61").Value = 1399.7407
$ws.Range("J61").Value = 4678.143
$ws.Range("K61").Value = 1399.7407
$ws.Range("L61").Value = 4678.143
$ws.Range("M61").Value = -1197.7407
$ws.Range("N61").Value = -5082.143

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 2519.195
$ws.Range("I113").Value = 1399.7407
$ws.Range("J113").Value = 4678.143
$ws.Range("K113").Value = 1399.7407
$ws.Range("L113").Value = 4678.143
$ws.Range("M113").Value = 770.2592999999999
$ws.Range("N113").Value = -9018.143

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3491.0334
$ws.Range("I122").Value = 3303.1155
$ws.Range("J122").Value = 4712.5
$ws.Range("K122").Value = 9909.3465
$ws.Range("L122").Value = 14137.5
$ws.Range("M122").Value = -7459.3465
$ws.Range("N122").Value = -19037.5

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 7010.5625
$ws.Range("I136").Value = 4183.8
$ws.Range("K136").Value = 12551.4
$ws.Range("M136").Value = -10001.4

$ws = $wb.Worksheets.Item("WVR")
# Row 49: A Leg Up on the Cold / Linen Tights
$ws.Range("H49").Value = 25997.1
$ws.Range("J49").Value = 26219.223
$ws.Range("L49").Value = 26219.223
$ws.Range("N49").Value = -26679.223

# Row 129: Lifetime of Gleaning / Scarlet Moko Beret of Gathering
$ws.Range("H129").Value = 40390
$ws.Range("I129").Value = 40390
$ws.Range("K129").Value = 40390
$ws.Range("M129").Value = -35390

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 2738.353
$ws.Range("I136").Value = 2131.2144
$ws.Range("K136").Value = 6393.6432
$ws.Range("M136").Value = -3843.6432
